# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# (commit: Updated cryptos list on Tue Sep 17 16:33:47 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume columns hold text-formatted numbers (thousand-dot separators,
# padded percentages, fixed decimal places). Force text format first so Excel
# does not reinterpret the assigned strings as numeric values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '60.885.13'
$ws.Range('E2').Value = '  +5.19%  '
$ws.Range('D3').Value = '2.370.93'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '546.89'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').Value = '133.41'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +2.70%  '
$ws.Range('D9').Value = '2.371.33'
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D11').Value = '5.52'
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('D14').Value = '24.23'
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').Value = '2.794.06'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('D16').Value = '60.772.00'
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +1.92%  '
$ws.Range('D18').Value = '2.319.41'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '10.75'
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').Value = '7.04'
$ws.Range('E20').Value = '  +10.26%  '
$ws.Range('D21').Value = '4.22'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').Value = '318.40'
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '63.77'
$ws.Range('E24').Value = '  +1.65%  '
$ws.Range('D25').Value = '0.173'
$ws.Range('E25').Value = '  +4.00%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '8.12'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('D28').Value = '1.35'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').Value = '1.75'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '171.90'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0749'
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('D32').Value = '1.14'
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('D33').Value = '5.92'
$ws.Range('E33').Value = '  +2.87%  '
$ws.Range('D34').Value = '1.41'
$ws.Range('E34').Value = '  +13.93%  '
$ws.Range('D35').Value = '0.385'
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('D36').Value = '18.11'
$ws.Range('E36').Value = '  +2.35%  '
$ws.Range('D38').Value = '4.25'
$ws.Range('E38').Value = '  +8.94%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '321.48'
$ws.Range('E40').Value = '  +11.21%  '
$ws.Range('D41').Value = '1.56'
$ws.Range('E41').Value = '  +4.92%  '
$ws.Range('D42').Value = '38.38'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').Value = '144.72'
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('D44').Value = '3.49'
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('D45').Value = '0.0955'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('D46').Value = '19.59'
$ws.Range('E46').Value = '  +8.16%  '
$ws.Range('D47').Value = '0.0503'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').Value = '0.567'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').Value = '0.0214'
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').Value = '0.0₆0203'
$ws.Range('E51').Value = '  +1.73%  '
